$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address -> new text value.
# NumberFormat is forced to text ("@") before assignment so that
# values which look numeric (e.g. "260.10", "0.9999") are stored
# verbatim as strings instead of being coerced/rounded by Excel,
# then the style is reset to "Normal" so no stray number format
# is left behind on the cell.
$updates = [ordered]@{
    "D2" = "26.638.96"
    "E2" = "  +0.83%  "
    "D3" = "1.845.05"
    "E3" = "  +0.21%  "
    "E4" = "  -0.03%  "
    "D5" = "260.10"
    "E5" = "  -0.39%  "
    "D6" = "0.9999"
    "E6" = "  -0.05%  "
    "D7" = "0.5268"
    "E7" = "  +1.80%  "
    "D8" = "0.3158"
    "E8" = "  -3.34%  "
    "D9" = "0.06799"
    "E9" = "  +0.34%  "
    "D10" = "18.78"
    "E10" = "  +0.45%  "
    "D11" = "0.7843"
    "E11" = "  +1.21%  "
    "D12" = "0.07767"
    "E12" = "  +0.28%  "
    "D13" = "1.849.33"
    "E13" = "  +0.49%  "
    "D14" = "88.09"
    "E14" = "  +0.31%  "
    "D15" = "5.018"
    "E15" = "  +0.38%  "
    "E16" = "  +0.01%  "
    "E17" = "  -0.38%  "
    "E18" = "  -0.03%  "
    "D19" = "0.000007927"
    "E19" = "  -0.05%  "
    "D20" = "26.669.32"
    "E20" = "  +0.80%  "
    "D21" = "2.078.14"
    "E21" = "  +0.00%  "
    "E22" = "  +0.02%  "
    "E23" = "  -0.28%  "
    "D24" = "9.319"
    "E24" = "  -2.33%  "
    "D25" = "2.219"
    "E25" = "  +1.02%  "
    "D26" = "142.40"
    "E26" = "  -2.26%  "
    "D27" = "1.678"
    "E27" = "  +1.62%  "
    "D28" = "17.03"
    "E28" = "  +0.29%  "
    "D29" = "111.04"
    "E29" = "  -0.56%  "
    "D30" = "4.197"
    "E31" = "  +0.29%  "
    "D32" = "4.083"
    "E32" = "  -0.97%  "
    "E33" = "  +1.36%  "
    "D34" = "0.7295"
    "E34" = "  +1.52%  "
    "D35" = "1.142"
    "D36" = "2.856"
    "E36" = "  +0.23%  "
    "D37" = "3.097"
    "E37" = "  +0.39%  "
    "D38" = "2.279"
    "E38" = "  +2.49%  "
    "D39" = "0.01732"
    "E39" = "  -2.47%  "
    "D40" = "0.4778"
    "E40" = "  -1.10%  "
    "D41" = "0.9002"
    "E41" = "  +0.15%  "
    "D42" = "109.89"
    "E42" = "  -1.85%  "
    "D43" = "5.952"
    "E43" = "  -2.10%  "
    "E44" = "  -0.03%  "
    "D45" = "7.692"
    "E45" = "  -0.45%  "
    "D46" = "0.4169"
    "E46" = "  +0.54%  "
    "D47" = "9.056"
    "E47" = "  +0.43%  "
    "D48" = "0.1238"
    "E48" = "  +1.70%  "
    "D49" = "0.05811"
    "E49" = "  -2.48%  "
    "D50" = "34.78"
    "E50" = "  -0.70%  "
    "D51" = "0.8923"
    "E51" = "  +0.75%  "
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}
